# Update the correlation table (macro_corr_educ_gr) with refreshed values
# exported from Stata. The table cells hold rounded-correlation strings
# (e.g. "-0.09", "0.0", "-0.2*") that must stay text, not be reinterpreted
# as numbers, so force each target cell to Text format before writing and
# then drop back to the sheet's default "Normal" style so no visible
# formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "B2" "-0.09"
Set-TextValue $ws "C2" "0.03"
Set-TextValue $ws "D2" "-0.09"
Set-TextValue $ws "C3" "0.02"
Set-TextValue $ws "D3" "-0.04"
Set-TextValue $ws "E3" "0.16"
Set-TextValue $ws "B4" "-0.04"
Set-TextValue $ws "D4" "0.1"
Set-TextValue $ws "E4" "0.11"
Set-TextValue $ws "G4" "-0.07"
Set-TextValue $ws "B5" "0.0"
Set-TextValue $ws "C5" "-0.06"
Set-TextValue $ws "C6" "-0.2*"
Set-TextValue $ws "D6" "-0.11"
Set-TextValue $ws "E6" "-0.08"
Set-TextValue $ws "F6" "-0.06"
Set-TextValue $ws "G6" "0.02"
Set-TextValue $ws "B7" "-0.08"
Set-TextValue $ws "C7" "-0.14"
Set-TextValue $ws "E7" "-0.2*"
Set-TextValue $ws "G7" "0.09"
Set-TextValue $ws "B8" "-0.04"
Set-TextValue $ws "E8" "-0.11"
Set-TextValue $ws "B9" "0.01"
Set-TextValue $ws "D9" "0.03"
Set-TextValue $ws "G9" "0.05"
Set-TextValue $ws "B10" "-0.21*"
Set-TextValue $ws "D10" "-0.13"
Set-TextValue $ws "F10" "-0.09"
Set-TextValue $ws "G10" "-0.02"
Set-TextValue $ws "C11" "0.04"
Set-TextValue $ws "D11" "-0.15"
Set-TextValue $ws "E11" "0.01"
Set-TextValue $ws "F11" "0.05"
Set-TextValue $ws "B12" "-0.15"
Set-TextValue $ws "C12" "-0.17"
Set-TextValue $ws "D12" "-0.08"
Set-TextValue $ws "E12" "0.0"
Set-TextValue $ws "F12" "0.03"
Set-TextValue $ws "B13" "0.08"
Set-TextValue $ws "D13" "0.1"
Set-TextValue $ws "E13" "-0.04"
